$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column (H) — copy the header formatting from the
# neighboring "sum" header cell (G1) so the new header matches the
# existing header style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the Save values for the two data rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
